$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.268.48"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.841.84"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.04"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6705"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.54%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07431"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2938"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.86"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07719"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.81%  "
$ws.Range("D12").Value = "1.826.84"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.007"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6713"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.71%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "85.89"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.148"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("D17").Value = "29.236.33"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008309"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.89"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.151"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.702"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  -3.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.04"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.511"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  -2.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.067"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.195"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05301"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.878"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7519"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.136"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.683"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").Value = "1.320.44"
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01805"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.727"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9209"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.962"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.08487"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +15.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.007"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.85"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("D45").Value = "1.975.15"
$ws.Range("E45").Value = "  -0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5164"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.775"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.86"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.140"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05943"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.09%  "
